# Update the cryptocurrency price/volume table to the new snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.783.41"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "'3.301.89"
$ws.Range("E3").Value = "  +1.13%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'558.00"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").Value = "'183.88"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'3.293.91"
$ws.Range("E8").Value = "  +1.08%  "

$ws.Range("E9").Value = "  -3.03%  "

$ws.Range("E10").Value = "  -5.25%  "

$ws.Range("D11").Value = "'0.573"
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").Value = "'45.43"
$ws.Range("E12").Value = "  -3.03%  "

$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").Value = "'3.831.65"
$ws.Range("E14").Value = "  +1.37%  "

$ws.Range("D15").Value = "'8.37"
$ws.Range("E15").Value = "  -2.52%  "

$ws.Range("D16").Value = "'574.22"
$ws.Range("E16").Value = "  -9.07%  "

$ws.Range("D17").Value = "'65.731.76"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").Value = "'3.300.01"
$ws.Range("E19").Value = "  +1.14%  "

$ws.Range("E20").Value = "  -2.69%  "

$ws.Range("D21").Value = "'10.80"
$ws.Range("E21").Value = "  -4.34%  "

$ws.Range("D22").Value = "'0.886"
$ws.Range("E22").Value = "  -1.56%  "

$ws.Range("D23").Value = "'17.60"
$ws.Range("E23").Value = "  -3.33%  "

$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("D25").Value = "'97.27"
$ws.Range("E25").Value = "  -8.64%  "

$ws.Range("D26").Value = "'3.91"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").Value = "'5.90"
$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("E28").Value = "  +0.40%  "

$ws.Range("E29").Value = "  -2.64%  "

$ws.Range("D30").Value = "'30.41"
$ws.Range("E30").Value = "  +1.01%  "

$ws.Range("D31").Value = "'8.38"
$ws.Range("E31").Value = "  -2.83%  "

$ws.Range("D32").Value = "'6.55"
$ws.Range("E32").Value = "  +5.53%  "

$ws.Range("D33").Value = "'3.66"
$ws.Range("E33").Value = "  -6.22%  "

$ws.Range("D34").Value = "'554.53"
$ws.Range("E34").Value = "  +6.57%  "

$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.102"
$ws.Range("E36").Value = "  -1.74%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "'3.715.29"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").Value = "'55.49"
$ws.Range("E39").Value = "  -3.91%  "

$ws.Range("D40").Value = "'32.96"
$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("E41").Value = "  -3.36%  "

$ws.Range("E42").Value = "  -7.17%  "

$ws.Range("D43").Value = "'0.0₃0677"
$ws.Range("E43").Value = "  -7.02%  "

$ws.Range("E44").Value = "  +3.60%  "

$ws.Range("D45").Value = "'2.55"
$ws.Range("E45").Value = "  -5.27%  "

$ws.Range("E46").Value = "  -1.64%  "

$ws.Range("D47").Value = "'0.0405"
$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("D48").Value = "'3.08"
$ws.Range("E48").Value = "  -8.96%  "

$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").Value = "'0.126"
$ws.Range("E50").Value = "  -2.37%  "

$ws.Range("D51").Value = "'2.49"
$ws.Range("E51").Value = "  -3.62%  "
